$wb = $excel.ActiveWorkbook

# Template sheet to copy cell formatting from (the "2020-2021" sheet,
# which has the same A-column layout: region names with two "section
# header" rows styled distinctly).
$template = $wb.Worksheets.Item("2020-2021")

# Add the new sheet as the last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2024-2025"

# Row 1: header labels (these reuse the existing shared strings).
$newSheet.Range("A1").Value = "Región"
$newSheet.Range("B1").Value = "Total"
$newSheet.Range("C1").Value = "A solicitud de la parte peticionaria"
$newSheet.Range("D1").Value = "Otra razón"

# Column A region labels, rows 2-16.
$regions = @{
    2  = "Aguadilla"
    3  = "Aibonito"
    4  = "Arecibo"
    5  = "Bayamón"
    6  = "Caguas"
    7  = "Carolina"
    8  = "Fajardo"
    9  = "Guayama"
    10 = "Humacao"
    11 = "Mayagüez"
    12 = "Ponce"
    13 = "San Juan"
    14 = "Utuado"
    15 = "No indica"
    16 = "Total"
}

foreach ($r in $regions.Keys) {
    $newSheet.Range("A$r").Value = $regions[$r]
}

# Rows 2 and 13 carry the "section header" style from the template
# (template rows A2 and A13 use that same distinct style), so copy the
# formatting over from there instead of re-deriving it.
[void]$template.Range("A2").Copy()
[void]$newSheet.Range("A2").PasteSpecial(-4122)

[void]$template.Range("A13").Copy()
[void]$newSheet.Range("A13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New data entered for this period.
$newSheet.Range("B13").Value = 1
$newSheet.Range("D13").Value = 1

# Update selections on the other sheets to match the final saved state.
$ws1 = $wb.Worksheets.Item("2020-2021")
[void]$ws1.Range("C16").Select()

$ws2 = $wb.Worksheets.Item("2021-2022")
[void]$ws2.Range("A10").Select()

$ws3 = $wb.Worksheets.Item("2022-2023")
[void]$ws3.Range("A2").Select()

$ws4 = $wb.Worksheets.Item("2023-2024")
[void]$ws4.Range("A1:D16").Select()

# Make the new sheet the active/selected tab, with the selection left
# where data entry ended.
[void]$newSheet.Activate()
[void]$newSheet.Range("D18").Select()
